$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.463.93"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "2.781.17"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.45"
$ws.Range("E5").Value = "  -2.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.29"
$ws.Range("E6").Value = "  -1.39%  "

$ws.Range("E7").Value = "  -1.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  +4.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.28"
$ws.Range("E10").Value = "  -2.31%  "

$ws.Range("E11").Value = "  +1.39%  "

$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.89"
$ws.Range("E13").Value = "  +1.84%  "

$ws.Range("E14").Value = "  +2.21%  "

$ws.Range("D15").Value = "3.217.00"
$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("D16").Value = "2.805.01"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("E17").Value = "  -1.47%  "

$ws.Range("D18").Value = "51.475.08"
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.69"
$ws.Range("E19").Value = "  +2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.10"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.33"
$ws.Range("E21").Value = "  +1.41%  "

$ws.Range("D22").Value = "0.0₃0966"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.58"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.27"
$ws.Range("E24").Value = "  -1.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  -0.29%  "

$ws.Range("E27").Value = "  -2.74%  "

$ws.Range("E28").Value = "  +2.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.25"
$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.18"
$ws.Range("E30").Value = "  +7.82%  "

$ws.Range("E31").Value = "  -2.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.26"
$ws.Range("E32").Value = "  +8.82%  "

$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("E34").Value = "  +8.10%  "

$ws.Range("E35").Value = "  -6.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0851"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.47"
$ws.Range("E38").Value = "  -3.13%  "

$ws.Range("E39").Value = "  -3.19%  "

$ws.Range("E40").Value = "  -1.95%  "

$ws.Range("E41").Value = "  -0.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.48"
$ws.Range("E42").Value = "  -5.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.00"
$ws.Range("E43").Value = "  +0.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.88"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("E45").Value = "  -2.70%  "

$ws.Range("D46").Value = "2.130.06"
$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.34"
$ws.Range("E47").Value = "  +2.56%  "

$ws.Range("E48").Value = "  +5.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.225"
$ws.Range("E49").Value = "  +17.83%  "

$ws.Range("E50").Value = "  -5.71%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.34"
$ws.Range("E51").Value = "  +8.34%  "

